$d = $word.ActiveDocument

$replacements = @(
    @{old="88×62=5456"; new="94×22=2068"},
    @{old="81×53=4293"; new="63×91=5733"},
    @{old="90×91=8190"; new="68×52=3536"},
    @{old="49×48=2352"; new="65×87=5655"},
    @{old="43×62=2666"; new="41×30=1230"},
    @{old="15×96=1440"; new="97×62=6014"},
    @{old="62×60=3720"; new="37×96=3552"},
    @{old="39×52=2028"; new="63×52=3276"},
    @{old="45×69=3105"; new="32×39=1248"},
    @{old="44×85=3740"; new="70×52=3640"},
    @{old="93×88=8184"; new="76×14=1064"},
    @{old="20×76=1520"; new="87×42=3654"},
    @{old="54×72=3888"; new="57×13=741"},
    @{old="74×62=4588"; new="24×47=1128"},
    @{old="92×88=8096"; new="18×64=1152"},
    @{old="53×43=2279"; new="70×75=5250"},
    @{old="58×82=4756"; new="57×34=1938"},
    @{old="19×62=1178"; new="86×51=4386"},
    @{old="47×47=2209"; new="53×42=2226"},
    @{old="79×75=5925"; new="83×68=5644"},
    @{old="16×31=496"; new="65×86=5590"},
    @{old="32×15=480"; new="35×78=2730"},
    @{old="82×98=8036"; new="54×95=5130"},
    @{old="85×20=1700"; new="50×87=4350"},
    @{old="27×37=999"; new="41×99=4059"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
